$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cells that introduce brand-new shared strings first, and in the
# same order they were first introduced in the authored workbook, so the
# shared-string table indices line up exactly.
$ws.Range("H1").Value = "MK1"
$ws.Range("H21").Value = "14,17.3"
$ws.Range("H22").Value = "21,23"
$ws.Range("H10").Value = "20,22.2"
$ws.Range("H14").Value = "9,12"
$ws.Range("H23").Value = "24.2,29.2"

# Remaining column H data values (rows 2-23, with rows 16 and 17 left blank)
# These all reuse strings/values that already exist in the workbook.
$ws.Range("H2").Value = "15,17"
$ws.Range("H3").Value = "12,13"
$ws.Range("H4").Value = "1,2"
$ws.Range("H5").Value = "6,9.3"
$ws.Range("H6").Value = "8,9"
$ws.Range("H7").Value = "9,11"
$ws.Range("H8").Value = "10,12"
$ws.Range("H9").Value = "16,18"
$ws.Range("H11").Value = 10
$ws.Range("H12").Value = 30.2
$ws.Range("H13").Value = "14,16"
$ws.Range("H15").Value = 11
$ws.Range("H18").Value = 15
$ws.Range("H19").Value = 15
$ws.Range("H20").Value = 14

# Match column H formatting (font size 12 + centered) to the same style used by columns B-G
$ws.Range("H1:H15").Font.Size = 12
$ws.Range("H1:H15").HorizontalAlignment = -4108
$ws.Range("H18:H23").Font.Size = 12
$ws.Range("H18:H23").HorizontalAlignment = -4108

# Update the active selection to reflect where editing left off
$null = $ws.Range("H26").Select()
